$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the rotation between rows 3, 4 and 5.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AO")

# Capture the current ("before") values for rows 3-5 in the affected columns.
# NOTE: use Value2 (not the parameterized Value property) to reliably read
# the underlying cell value via COM interop.
$before = @{}
foreach ($r in 3..5) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# The rotation: new row3 <- old row5, new row4 <- old row3, new row5 <- old row4.
$mapping = @{ 3 = 5; 4 = 3; 5 = 4 }

foreach ($r in 3..5) {
    $srcRow = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $before[$srcRow][$c]
    }
}
